# "Edit item" feature: rename a duplicated "glop" item and add a brand-new
# "lop" item to the Inventory sheet, plus record the new sales that were
# logged for "glop" on the Sales sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Inventory sheet
# ---------------------------------------------------------------
$inv = $wb.Worksheets.Item("Inventory")

# Header tweak: "45% Profit" -> "Profit"
$inv.Cells.Item(1, 9).Value = "Profit"

# Insert a brand new item "lop" as the new row 2 (pushes existing rows down)
$inv.Rows.Item(2).Insert()

$inv.Cells.Item(2, 1).Value  = "lop"
$inv.Cells.Item(2, 2).Value  = "stock"
$inv.Cells.Item(2, 3).Value  = "chn"
$inv.Cells.Item(2, 4).Value  = "'12"
$inv.Cells.Item(2, 5).Value  = "pcs"
$inv.Cells.Item(2, 6).Value  = 23
$inv.Cells.Item(2, 7).Value  = 921321546465
$inv.Cells.Item(2, 8).Value  = "'10.5"
$inv.Cells.Item(2, 9).Value  = "'4.73"
$inv.Cells.Item(2, 10).Value = "'15.23"
$inv.Cells.Item(2, 11).Value = "'0.71"
$inv.Cells.Item(2, 12).Value = "'15.93"
$inv.Cells.Item(2, 13).Value = 23
$inv.Cells.Item(2, 14).Value = "'0.00"
$inv.Cells.Item(2, 15).Value = "'0.00"

# Rename the item that is now on row 4 (was the duplicate "glop" row)
$inv.Cells.Item(4, 1).Value = "glplpp"

# ---------------------------------------------------------------
# Sales sheet: log the new sales made for "glop"
# ---------------------------------------------------------------
$sales = $wb.Worksheets.Item("Sales")

$sales.Cells.Item(4, 1).Value = "glop"
$sales.Cells.Item(4, 2).Value = 10.5
$sales.Cells.Item(4, 3).Value = 10
$sales.Cells.Item(4, 4).Value = 105
$sales.Cells.Item(4, 5).Value = "lewi"
$sales.Cells.Item(4, 6).Value = "2023-01-03 12:32:11"

$sales.Cells.Item(5, 1).Value = "glop"
$sales.Cells.Item(5, 2).Value = 10.5
$sales.Cells.Item(5, 3).Value = 10
$sales.Cells.Item(5, 4).Value = 105
$sales.Cells.Item(5, 5).Value = "lewi"
$sales.Cells.Item(5, 6).Value = "2023-01-03 12:33:53"

$sales.Cells.Item(6, 1).Value = "glop"
$sales.Cells.Item(6, 2).Value = 10.5
$sales.Cells.Item(6, 3).Value = 10
$sales.Cells.Item(6, 4).Value = 105
$sales.Cells.Item(6, 5).Value = "lewi"
$sales.Cells.Item(6, 6).Value = "2023-01-03 12:36:35"
